# Append two new job listings to the top of the "ランサーズ" sheet's table.
# New row 2: n8n job posting (2025-09-18 18:25:58)
# New row 5 (after the old-row-3 content that lands at row 4): 1688アリババ job posting
# Every other existing row shifts down (2 rows pushed into new slot 2, 2 more rows
# pushed into new slot 5) and keeps its original data/formatting untouched.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- 1. Make room: insert a blank row at row 2, then another blank row at row 5. ---
# (Row data/number formats shift down automatically; Hyperlinks do NOT follow the
# shift in this engine, so we rebuild the whole Hyperlinks collection afterwards.)
$ws.Rows.Item(2).Insert()
$ws.Rows.Item(5).Insert()

# --- 2. Fill the new row 2 ---
$ws.Range("A2").Value = "2025-09-18 18:25:58"
$ws.Range("B2").Value = "【報酬4万円〜|相談可能】n8n構築者募集|AIワークフロー構築が得意な方を探しています"
$ws.Range("C2").Value = "システム開発"
$ws.Range("D2").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E2").Value = "期限情報なし"
$ws.Range("F2").Value = "https://www.lancers.jp/work/detail/5396220"
$ws.Range("F2").Style = "Hyperlink"
$ws.Range("G2").Value = 303
$ws.Range("H2").Value = "🔥AI,Ai"

# --- 3. Fill the new row 5 ---
$ws.Range("A5").Value = "2025-09-18 18:25:58"
$ws.Range("B5").Value = "1688アリババの商品情報の抽出のスクレイピングの開発 exe形式の自動ツール"
$ws.Range("C5").Value = "システム開発"
$ws.Range("D5").Value = "50,000 円 ~ 100,000 円 / 固定"
$ws.Range("E5").Value = "期限情報なし"
$ws.Range("F5").Value = "https://www.lancers.jp/work/detail/5387065"
$ws.Range("F5").Style = "Hyperlink"
$ws.Range("G5").Value = 173
$ws.Range("H5").Value = "◆ツール,開発"

# --- 4. Column B is now wider (40 -> 46 stored character-width units). ---
# Excel's ColumnWidth getter/setter is offset from the stored <col width="..">
# value by the standard 5/6 character padding, so subtract that to land on an
# exact stored width of 46.
$ws.Columns.Item(2).ColumnWidth = 46 - 5/6

# --- 5. Rebuild hyperlinks F2:F19 so every URL cell is a live hyperlink again. ---
# (Inserting rows does not carry the Hyperlinks collection along, so wipe and
# re-add in row order; this also matches how the source file is regenerated.)
$ws.Range("F2").Hyperlinks.Delete()

for ($r = 2; $r -le 19; $r++) {
    $cell = $ws.Cells.Item($r, 6)
    $url = $cell.Value
    $ws.Hyperlinks.Add($cell, $url)
}

Write-Output "done"
